$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted as row 14; every subsequent record
# (previously rows 14-119) shifts down by one row, ending at row 120.
$ws.Rows.Item(14).Insert()

$ws.Range("A14").Value = 5
$ws.Range("B14").Value = "Macroferia Regional de Talca"
$ws.Range("C14").Value = "Maule"
$ws.Range("D14").Value = 45168
$ws.Range("D14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E14").Value = 7
$ws.Range("F14").Value = 100112040
$ws.Range("G14").Value = "Cilantro"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 150
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 10000
$ws.Range("N14").Value = "`$/caja 36 atados"
$ws.Range("O14").Value = "Región Metropolitana"
$ws.Range("P14").Value = 278
$ws.Range("Q14").Value = 36
$ws.Range("R14").Value = "Hortaliza"
